$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# New test row data references resume headline paths. The shared strings
# for these values are re-saved in "H16 first, then H10, H11, H12" order
# to match the regenerated shared string table ordering, and a leading
# apostrophe preserves the existing quote-prefix cell style (these text
# values look like they could start a formula/relative path).
$ws.Range("H16").Value = "'D:\\Naresh_Resume.pdf"
$ws.Range("H10").Value = "'D:\\K_Thrinath.docx"
$ws.Range("H11").Value = "'D:\\Sandeep_Resume.pdf"
$ws.Range("H12").Value = "'D:\\Uday_Resume.docx"
